$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every written cell to plain text ("@") so Excel does not
# reinterpret values like "0.998", "1.12" or "63.356.38" as numbers,
# then reset the style back to "Normal" so no stray cell style (s=...)
# is left behind, matching the original (unstyled) data cells.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '63.356.38'
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  +0.53%  '
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.028.65'
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -2.94%  '
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '558.63'
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  -0.20%  '
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '156.60'
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -4.10%  '
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  -0.12%  '
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  -4.99%  '
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '3.038.82'
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  -2.11%  '
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  -2.49%  '
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '6.42'
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -4.44%  '
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  -3.15%  '
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '3.556.07'
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  -2.60%  '
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  -3.22%  '
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '63.317.38'
$cell.Style = "Normal"

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  +0.33%  '
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  -2.24%  '
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.025.19'
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  -2.60%  '
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  -2.22%  '
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '399.67'
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  -1.58%  '
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '5.12'
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  -1.01%  '
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '12.08'
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -3.58%  '
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.72'
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  -4.83%  '
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  +0.24%  '
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '65.45'
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  -4.66%  '
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  -5.03%  '
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.467'
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  -4.40%  '
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0979'
$cell.Style = "Normal"

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  -3.87%  '
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '8.73'
$cell.Style = "Normal"

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  +0.35%  '
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.994'
$cell.Style = "Normal"

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  -0.53%  '
$cell.Style = "Normal"

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  +0.02%  '
$cell.Style = "Normal"

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  -3.67%  '
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '162.60'
$cell.Style = "Normal"

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  +5.08%  '
$cell.Style = "Normal"

$cell = $ws.Range("B34")
$cell.NumberFormat = "@"
$cell.Value = 'NEARProtocol'
$cell.Style = "Normal"

$cell = $ws.Range("C34")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '4.79'
$cell.Style = "Normal"

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  -1.51%  '
$cell.Style = "Normal"

$cell = $ws.Range("B35")
$cell.NumberFormat = "@"
$cell.Value = 'Fetch.AI'
$cell.Style = "Normal"

$cell = $ws.Range("C35")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.12'
$cell.Style = "Normal"

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  +2.06%  '
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.06'
$cell.Style = "Normal"

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -2.42%  '
$cell.Style = "Normal"

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  -0.53%  '
$cell.Style = "Normal"

$cell = $ws.Range("B38")
$cell.NumberFormat = "@"
$cell.Value = 'Maker'
$cell.Style = "Normal"

$cell = $ws.Range("C38")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.553.78'
$cell.Style = "Normal"

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  -6.25%  '
$cell.Style = "Normal"

$cell = $ws.Range("B39")
$cell.NumberFormat = "@"
$cell.Value = 'Stacks'
$cell.Style = "Normal"

$cell = $ws.Range("C39")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.61'
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  -2.61%  '
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '22.87'
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  -3.49%  '
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -4.01%  '
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '37.73'
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  -1.10%  '
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.670'
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  -3.48%  '
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0605'
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  +0.04%  '
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.0251'
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -1.86%  '
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '5.14'
$cell.Style = "Normal"

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  -2.12%  '
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  -0.14%  '
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '20.62'
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  -1.02%  '
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '273.71'
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  -3.62%  '
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0947'
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  -2.40%  '
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '10.48'
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  +0.02%  '
$cell.Style = "Normal"

